# Generate Report for Handoff
# Replace the old GUID-named source markdown file references with the new
# GUID (032ceea2-d288-4b97-9cb2-4850214ad6e1) across the Overview / zh-cn /
# de-de sheets, refresh the cached handoff/xliff file names (new content
# hash), and bump the "generated at" / "handoff at" timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "6a43499c-4e72-4ddd-8a7b-e547cc1eaccd"
$newGuid = "032ceea2-d288-4b97-9cb2-4850214ad6e1"
$oldHash = "b23022cde0df3606918046fd60d27eeb28c47a00"
$newHash = "c921fb461b9a0234035f6a9bdaa9825c3ac0fcac"

$newFileName = "$newGuid.md"
$newPath     = "e2e\$newGuid.md"

# The external hyperlinks all point at the same GitHub blob URL (keyed off
# the old GUID); that underlying address is left as-is, only the cached
# display text is refreshed to show the new file name.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5172782d962a8c954ea735f7a0823bc0ebb426a/e2e/$oldGuid.md"

function Set-HyperlinkDisplay {
    param($ws, $cellAddr, $displayText)

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $displayText)
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPath
Set-HyperlinkDisplay $wsOverview "B2" $newPath
$wsOverview.Range("G2").Value = "2016-08-21 23:05:14"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newFileName
Set-HyperlinkDisplay $wsZhCn "A2" $newFileName
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 23:05:10"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newFileName
Set-HyperlinkDisplay $wsDeDe "A2" $newFileName
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-21 23:05:14"
